$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44231
$ws.Range("J2").Value = 120
$ws.Range("D3").Value = 44294
$ws.Range("J3").Value = 180
$ws.Range("D4").Value = 44165
$ws.Range("J4").Value = 68
$ws.Range("K4").Value = 2000
$ws.Range("L4").Value = 2000
$ws.Range("M4").Value = 2000
$ws.Range("P4").Value = 667
$ws.Range("D5").Value = 44488
$ws.Range("D6").Value = 44285
$ws.Range("J6").Value = 160
$ws.Range("D7").Value = 44428
$ws.Range("J7").Value = 160
$ws.Range("D8").Value = 44511
$ws.Range("D9").Value = 44274
$ws.Range("D10").Value = 44371
$ws.Range("J10").Value = 180
$ws.Range("D11").Value = 44503
$ws.Range("J11").Value = 160
$ws.Range("D12").Value = 44263
$ws.Range("J12").Value = 180
$ws.Range("D13").Value = 44410
$ws.Range("D14").Value = 44196
$ws.Range("J14").Value = 180
$ws.Range("D15").Value = 44413
$ws.Range("D16").Value = 44193
$ws.Range("J16").Value = 120
$ws.Range("K16").Value = 1800
$ws.Range("L16").Value = 1800
$ws.Range("M16").Value = 1800
$ws.Range("P16").Value = 600
$ws.Range("D17").Value = 44490
$ws.Range("D18").Value = 44351
$ws.Range("D19").Value = 44258
$ws.Range("J19").Value = 230
$ws.Range("D20").Value = 44343
$ws.Range("J20").Value = 180
$ws.Range("D21").Value = 44160
$ws.Range("J21").Value = 230
$ws.Range("D22").Value = 44365
$ws.Range("J22").Value = 180
$ws.Range("D23").Value = 44476
$ws.Range("D24").Value = 44411
$ws.Range("J24").Value = 120
$ws.Range("D25").Value = 44459
$ws.Range("J25").Value = 160
$ws.Range("D26").Value = 44281
$ws.Range("J26").Value = 160
$ws.Range("D27").Value = 44364
$ws.Range("D28").Value = 44313
$ws.Range("J28").Value = 130
$ws.Range("D29").Value = 44334
$ws.Range("J29").Value = 190
$ws.Range("D30").Value = 44319
$ws.Range("J30").Value = 190
$ws.Range("D31").Value = 44414
$ws.Range("J31").Value = 160
$ws.Range("D32").Value = 44316
$ws.Range("J32").Value = 160
$ws.Range("D33").Value = 44461
$ws.Range("J33").Value = 160
$ws.Range("D34").Value = 44466
$ws.Range("J34").Value = 160
$ws.Range("D35").Value = 44385
$ws.Range("J35").Value = 180
$ws.Range("D36").Value = 44379
$ws.Range("D37").Value = 44253
$ws.Range("K37").Value = 1500
$ws.Range("L37").Value = 1500
$ws.Range("M37").Value = 1500
$ws.Range("P37").Value = 500
$ws.Range("D38").Value = 44370
$ws.Range("J38").Value = 180
$ws.Range("D39").Value = 44484
$ws.Range("J39").Value = 160
$ws.Range("D40").Value = 44418
$ws.Range("J40").Value = 150
$ws.Range("D41").Value = 44427
$ws.Range("J41").Value = 160
$ws.Range("D42").Value = 44179
$ws.Range("J42").Value = 48
$ws.Range("K42").Value = 2000
$ws.Range("L42").Value = 2000
$ws.Range("M42").Value = 2000
$ws.Range("P42").Value = 667
$ws.Range("D43").Value = 44455
$ws.Range("D44").Value = 44280
$ws.Range("D45").Value = 44377
$ws.Range("J45").Value = 160
$ws.Range("D46").Value = 44369
$ws.Range("J46").Value = 180
$ws.Range("D47").Value = 44215
$ws.Range("J47").Value = 130
$ws.Range("D48").Value = 44397
$ws.Range("D49").Value = 44344
$ws.Range("J49").Value = 160
$ws.Range("D50").Value = 44188
$ws.Range("J50").Value = 180
$ws.Range("D51").Value = 44383
$ws.Range("D52").Value = 44267
$ws.Range("D53").Value = 44314
$ws.Range("D54").Value = 44412
$ws.Range("J54").Value = 160
$ws.Range("D55").Value = 44266
$ws.Range("J55").Value = 120
$ws.Range("D56").Value = 44473
$ws.Range("J56").Value = 160
$ws.Range("D57").Value = 44433
$ws.Range("J57").Value = 180
$ws.Range("D58").Value = 44159
$ws.Range("J58").Value = 120
$ws.Range("D59").Value = 44512
$ws.Range("D60").Value = 44286
$ws.Range("J60").Value = 160
$ws.Range("K60").Value = 1500
$ws.Range("L60").Value = 1500
$ws.Range("M60").Value = 1500
$ws.Range("P60").Value = 500
$ws.Range("D61").Value = 44335
$ws.Range("J61").Value = 160
$ws.Range("D62").Value = 44186
$ws.Range("J62").Value = 180
$ws.Range("D63").Value = 44460
$ws.Range("D64").Value = 44438
$ws.Range("D65").Value = 44392
$ws.Range("D66").Value = 44355
$ws.Range("J66").Value = 180
$ws.Range("D67").Value = 44489
$ws.Range("J67").Value = 160
$ws.Range("D68").Value = 44434
$ws.Range("J68").Value = 140
$ws.Range("D69").Value = 44497
$ws.Range("D70").Value = 44358
$ws.Range("D71").Value = 44399
$ws.Range("J71").Value = 120
$ws.Range("D72").Value = 44298
$ws.Range("D73").Value = 44482
$ws.Range("J73").Value = 160
$ws.Range("D74").Value = 44405
$ws.Range("J74").Value = 160
$ws.Range("D75").Value = 44250
$ws.Range("D76").Value = 44218
$ws.Range("J76").Value = 130
$ws.Range("D77").Value = 44273
$ws.Range("J77").Value = 160
$ws.Range("D78").Value = 44386
$ws.Range("J78").Value = 160
$ws.Range("D79").Value = 44435
$ws.Range("J79").Value = 810
$ws.Range("D80").Value = 44328
$ws.Range("J80").Value = 160
$ws.Range("D81").Value = 44277
$ws.Range("J81").Value = 160
$ws.Range("D82").Value = 44442
$ws.Range("J82").Value = 180
$ws.Range("D83").Value = 44516
$ws.Range("J83").Value = 150
$ws.Range("D84").Value = 44175
$ws.Range("J84").Value = 120
$ws.Range("D85").Value = 44168
$ws.Range("D86").Value = 44203
$ws.Range("D87").Value = 44475
$ws.Range("D88").Value = 44483
$ws.Range("J88").Value = 180
$ws.Range("D89").Value = 44217
$ws.Range("J89").Value = 120
$ws.Range("D90").Value = 44235
$ws.Range("J90").Value = 160
$ws.Range("D91").Value = 44200
$ws.Range("J91").Value = 120
$ws.Range("D92").Value = 44419
$ws.Range("J92").Value = 130
$ws.Range("D93").Value = 44162
$ws.Range("D94").Value = 44357
$ws.Range("J94").Value = 160
$ws.Range("D95").Value = 44244
$ws.Range("J95").Value = 110
$ws.Range("D96").Value = 44202
$ws.Range("J96").Value = 120
$ws.Range("D97").Value = 44333
$ws.Range("D98").Value = 44320
$ws.Range("D99").Value = 44252
$ws.Range("D100").Value = 44467
$ws.Range("J100").Value = 160
$ws.Range("D101").Value = 44264
$ws.Range("J101").Value = 120
$ws.Range("D102").Value = 44214
$ws.Range("J102").Value = 110
$ws.Range("D103").Value = 44167
$ws.Range("J103").Value = 150
$ws.Range("D104").Value = 44291
$ws.Range("J104").Value = 89
$ws.Range("K104").Value = 1800
$ws.Range("L104").Value = 1800
$ws.Range("M104").Value = 1800
$ws.Range("P104").Value = 600
$ws.Range("D105").Value = 44174
$ws.Range("J105").Value = 180
$ws.Range("D106").Value = 44293
$ws.Range("D107").Value = 44496
$ws.Range("J107").Value = 150
$ws.Range("D108").Value = 44326
$ws.Range("J108").Value = 120
$ws.Range("D109").Value = 44302
$ws.Range("J109").Value = 130
$ws.Range("D110").Value = 44308
$ws.Range("J110").Value = 160
$ws.Range("D111").Value = 44498
$ws.Range("J111").Value = 160
$ws.Range("D112").Value = 44420
$ws.Range("J112").Value = 160
$ws.Range("D113").Value = 44398
$ws.Range("J113").Value = 160
$ws.Range("D114").Value = 44396
$ws.Range("D115").Value = 44321
$ws.Range("J115").Value = 130
$ws.Range("D116").Value = 44208
$ws.Range("J116").Value = 160
$ws.Range("D117").Value = 44349
$ws.Range("D118").Value = 44477
$ws.Range("J118").Value = 160
$ws.Range("K118").Value = 1500
$ws.Range("L118").Value = 1500
$ws.Range("M118").Value = 1500
$ws.Range("P118").Value = 500
$ws.Range("D119").Value = 44487
$ws.Range("J119").Value = 160
$ws.Range("D120").Value = 44452
$ws.Range("J120").Value = 190
$ws.Range("D121").Value = 44211
$ws.Range("J121").Value = 120
$ws.Range("D122").Value = 44505
$ws.Range("J122").Value = 120
$ws.Range("D123").Value = 44204
$ws.Range("J123").Value = 180
$ws.Range("D124").Value = 44306
$ws.Range("D125").Value = 44509
$ws.Range("D126").Value = 44454
$ws.Range("J126").Value = 160
$ws.Range("D127").Value = 44189
$ws.Range("D128").Value = 44278
$ws.Range("J128").Value = 130
$ws.Range("D129").Value = 44265
$ws.Range("J129").Value = 120
$ws.Range("D130").Value = 44494
$ws.Range("J130").Value = 190
$ws.Range("D131").Value = 44300
$ws.Range("J131").Value = 160
$ws.Range("D132").Value = 44209
$ws.Range("D133").Value = 44237
$ws.Range("J133").Value = 130
$ws.Range("D134").Value = 44356
$ws.Range("J134").Value = 160
$ws.Range("D135").Value = 44469
$ws.Range("D136").Value = 44453
$ws.Range("J136").Value = 130
$ws.Range("D137").Value = 44446
$ws.Range("J137").Value = 180
$ws.Range("D138").Value = 44463
$ws.Range("D139").Value = 44245
$ws.Range("D140").Value = 44323
$ws.Range("D141").Value = 44229
$ws.Range("D142").Value = 44417
$ws.Range("J142").Value = 160
$ws.Range("D143").Value = 44445
$ws.Range("J143").Value = 180
$ws.Range("D144").Value = 44249
$ws.Range("D145").Value = 44342
$ws.Range("J145").Value = 260
$ws.Range("D146").Value = 44259
$ws.Range("J146").Value = 120
$ws.Range("D147").Value = 44216
$ws.Range("J147").Value = 80
$ws.Range("D148").Value = 44406
$ws.Range("D149").Value = 44295
$ws.Range("J149").Value = 120
$ws.Range("D150").Value = 44270
$ws.Range("J150").Value = 120
$ws.Range("D151").Value = 44363
$ws.Range("J151").Value = 130
$ws.Range("D152").Value = 44299
$ws.Range("J152").Value = 130
$ws.Range("D153").Value = 44257
$ws.Range("D154").Value = 44336
$ws.Range("J154").Value = 160
$ws.Range("D155").Value = 44372
$ws.Range("D156").Value = 44403
$ws.Range("J156").Value = 180
$ws.Range("D157").Value = 44195
$ws.Range("J157").Value = 180
$ws.Range("D158").Value = 44376
$ws.Range("D159").Value = 44474
$ws.Range("D160").Value = 44172
$ws.Range("J160").Value = 110
$ws.Range("D161").Value = 44421
$ws.Range("J161").Value = 180
$ws.Range("D162").Value = 44431
$ws.Range("J162").Value = 180
$ws.Range("D163").Value = 44239
$ws.Range("J163").Value = 120
$ws.Range("D164").Value = 44426
$ws.Range("D165").Value = 44448
$ws.Range("J165").Value = 160
$ws.Range("D166").Value = 44362
$ws.Range("J166").Value = 180
$ws.Range("D167").Value = 44210
$ws.Range("J167").Value = 120
$ws.Range("D168").Value = 44176
$ws.Range("J168").Value = 80
$ws.Range("D169").Value = 44301
$ws.Range("J169").Value = 130
$ws.Range("D170").Value = 44407
$ws.Range("J170").Value = 160
$ws.Range("D171").Value = 44284
$ws.Range("J171").Value = 180
$ws.Range("D172").Value = 44441
$ws.Range("J172").Value = 190
$ws.Range("D173").Value = 44279
$ws.Range("D174").Value = 44341
$ws.Range("J174").Value = 160
$ws.Range("D176").Value = 44350
$ws.Range("D177").Value = 44312
$ws.Range("D178").Value = 44382
$ws.Range("J178").Value = 160
$ws.Range("D179").Value = 44384
$ws.Range("D180").Value = 44329
$ws.Range("D181").Value = 44246
$ws.Range("J181").Value = 160
$ws.Range("D182").Value = 44491
$ws.Range("J182").Value = 160
$ws.Range("D183").Value = 44272
$ws.Range("D184").Value = 44305
$ws.Range("J184").Value = 180
$ws.Range("D185").Value = 44447
$ws.Range("D186").Value = 44425
$ws.Range("D187").Value = 44315
$ws.Range("J187").Value = 130
$ws.Range("D188").Value = 44348
$ws.Range("D189").Value = 44322
$ws.Range("J189").Value = 130
$ws.Range("D190").Value = 44495
$ws.Range("J190").Value = 160
$ws.Range("D191").Value = 44232
$ws.Range("J191").Value = 120
$ws.Range("D192").Value = 44327
$ws.Range("J192").Value = 190
$ws.Range("D193").Value = 44510
$ws.Range("J193").Value = 160
$ws.Range("D194").Value = 44161
$ws.Range("J194").Value = 180
$ws.Range("D195").Value = 44468
$ws.Range("J195").Value = 180
$ws.Range("D196").Value = 44517
$ws.Range("D197").Value = 44238
$ws.Range("J197").Value = 130
$ws.Range("D198").Value = 44391
$ws.Range("J198").Value = 160
$ws.Range("D199").Value = 44236
$ws.Range("D200").Value = 44251
$ws.Range("J200").Value = 80
$ws.Range("D201").Value = 44515
$ws.Range("D202").Value = 44330
$ws.Range("J202").Value = 160
$ws.Range("D203").Value = 44432
$ws.Range("J203").Value = 150
$ws.Range("D204").Value = 44181
$ws.Range("J204").Value = 90
$ws.Range("D205").Value = 44194
$ws.Range("J205").Value = 80
$ws.Range("D206").Value = 44271
$ws.Range("J206").Value = 180
$ws.Range("D207").Value = 44307
$ws.Range("J207").Value = 130
$ws.Range("D208").Value = 44400
$ws.Range("J208").Value = 160
$ws.Range("D209").Value = 44309
$ws.Range("J209").Value = 160
$ws.Range("K209").Value = 1500
$ws.Range("L209").Value = 1500
$ws.Range("M209").Value = 1500
$ws.Range("P209").Value = 500
$ws.Range("D210").Value = 44508
$ws.Range("A211").Value = 3
$ws.Range("B211").Value = "Femacal de La Calera"
$ws.Range("C211").Value = "Coquimbo"
$ws.Range("D211").Value = 44201
$ws.Range("D211").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E211").Value = 5
$ws.Range("F211").Value = 100112039
$ws.Range("G211").Value = "Ciboulette"
$ws.Range("H211").Value = "Sin especificar"
$ws.Range("I211").Value = "Primera"
$ws.Range("J211").Value = 120
$ws.Range("K211").Value = 1500
$ws.Range("L211").Value = 1500
$ws.Range("M211").Value = 1500
$ws.Range("N211").Value = '$/docena de atados'
$ws.Range("O211").Value = "Provincia de Quillota"
$ws.Range("P211").Value = 500
$ws.Range("Q211").Value = 3
$ws.Range("R211").Value = "Hortaliza"
